$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")

# Row 63: add Temps [h] = 6 and Travail effectué text (new shared string)
$ws.Range("C63").Value = 6
$ws.Range("D63").Value = "Finalisation du refactor, mise à jour en live des écrans"

# Row 64: add Date + Type (same as row 63's type "Implémentation")
$ws.Range("A64").Value = (Get-Date -Year 2023 -Month 6 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B64").Value = "Implémentation"

# Update sheet view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 53
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B64").Select() | Out-Null
$ws.Range("B66").Select() | Out-Null
